# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain the same underlying data set, and the same set of
# rows were updated with new counts in this commit.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 44
    3  = 368
    8  = 32
    12 = 1130
    13 = 1462
    17 = 97
    19 = 56
    20 = 98
    21 = 255
    23 = 302
    24 = 1676
    28 = 631
    29 = 302
    30 = 59
    31 = 3961
    33 = 459
    34 = 233
    35 = 1008
    36 = 100
    37 = 45
    39 = 100
    40 = 46
    41 = 12
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
